$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.242.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.112.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.478"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.628.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.202.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.113.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "491.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -6.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "391.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.802.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  -8.18%  "
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  -2.73%  "
